$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3.0 ; $ws.Range("G2").Value = 51.10859366666667 ; $ws.Range("H2").Value = 153.325781 ; $ws.Range("I2").Value = 0.02866730932766026 ; $ws.Range("J2").Value = 0.02866730932766026 ; $ws.Range("K2").Value = 3.0 ; $ws.Range("M2").Value = 4.752338666666668 ; $ws.Range("N2").Value = 14.257016 ; $ws.Range("O2").Value = 0.04151639666945049 ; $ws.Range("P2").Value = 0.04151639666945049 ; $ws.Range("Q2").Value = 242.8853458810552 ; $ws.Range("R2").Value = 2185.968112929496 ; $ws.Range("S2").Value = 0.001190163385492981 ; $ws.Range("T2").Value = 0.001190163385492981
$ws.Range("E3").Value = 3.0 ; $ws.Range("G3").Value = 51.10859366666667 ; $ws.Range("H3").Value = 153.325781 ; $ws.Range("I3").Value = 0.02866730932766026 ; $ws.Range("J3").Value = 0.02866730932766026 ; $ws.Range("K3").Value = 3.0 ; $ws.Range("M3").Value = 7.095953333333333 ; $ws.Range("N3").Value = 21.28786 ; $ws.Range("O3").Value = 0.06199019766855336 ; $ws.Range("P3").Value = 0.06199019766855336 ; $ws.Range("Q3").Value = 362.6641955909623 ; $ws.Range("R3").Value = 3263.97776031866 ; $ws.Range("S3").Value = 0.001777092171847223 ; $ws.Range("T3").Value = 0.001777092171847223
$ws.Range("E4").Value = 3.0 ; $ws.Range("G4").Value = 51.10859366666667 ; $ws.Range("H4").Value = 153.325781 ; $ws.Range("I4").Value = 0.02866730932766026 ; $ws.Range("J4").Value = 0.02866730932766026 ; $ws.Range("K4").Value = 3.0 ; $ws.Range("M4").Value = 27.34521433333333 ; $ws.Range("N4").Value = 82.03564300000001 ; $ws.Range("O4").Value = 0.238887597223811 ; $ws.Range("P4").Value = 0.238887597223811 ; $ws.Range("Q4").Value = 1397.575448090243 ; $ws.Range("R4").Value = 12578.17903281218 ; $ws.Range("S4").Value = 0.006848264644156502 ; $ws.Range("T4").Value = 0.006848264644156502
$ws.Range("E5").Value = 3.0 ; $ws.Range("G5").Value = 51.10859366666667 ; $ws.Range("H5").Value = 153.325781 ; $ws.Range("I5").Value = 0.02866730932766026 ; $ws.Range("J5").Value = 0.02866730932766026 ; $ws.Range("K5").Value = 3.0 ; $ws.Range("M5").Value = 19.73820233333333 ; $ws.Range("N5").Value = 59.214607 ; $ws.Range("O5").Value = 0.1724327946912327 ; $ws.Range("P5").Value = 0.1724327946912327 ; $ws.Range("Q5").Value = 1008.791762764785 ; $ws.Range("R5").Value = 9079.125864883068 ; $ws.Range("S5").Value = 0.004943184263646499 ; $ws.Range("T5").Value = 0.004943184263646499
$ws.Range("E6").Value = 3.0 ; $ws.Range("G6").Value = 51.10859366666667 ; $ws.Range("H6").Value = 153.325781 ; $ws.Range("I6").Value = 0.02866730932766026 ; $ws.Range("J6").Value = 0.02866730932766026 ; $ws.Range("K6").Value = 3.0 ; $ws.Range("M6").Value = 51.532109 ; $ws.Range("N6").Value = 154.596327 ; $ws.Range("O6").Value = 0.4501841363839443 ; $ws.Range("P6").Value = 0.4501841363839443 ; $ws.Range("Q6").Value = 2633.733619667376 ; $ws.Range("R6").Value = 23703.60257700639 ; $ws.Range("S6").Value = 0.01290556789212412 ; $ws.Range("T6").Value = 0.01290556789212412
$ws.Range("E7").Value = 3.0 ; $ws.Range("G7").Value = 51.10859366666667 ; $ws.Range("H7").Value = 153.325781 ; $ws.Range("I7").Value = 0.02866730932766026 ; $ws.Range("J7").Value = 0.02866730932766026 ; $ws.Range("K7").Value = 3.0 ; $ws.Range("M7").Value = 4.005140333333333 ; $ws.Range("N7").Value = 12.015421 ; $ws.Range("O7").Value = 0.03498887736300817 ; $ws.Range("P7").Value = 0.03498887736300817 ; $ws.Range("Q7").Value = 204.6970898743112 ; $ws.Range("R7").Value = 1842.273808868801 ; $ws.Range("S7").Value = 0.001003036970392925 ; $ws.Range("T7").Value = 0.001003036970392925
$ws.Range("E8").Value = 3.0 ; $ws.Range("G8").Value = 1523.209716666667 ; $ws.Range("H8").Value = 4569.62915 ; $ws.Range("I8").Value = 0.854383206146807 ; $ws.Range("J8").Value = 0.854383206146807 ; $ws.Range("K8").Value = 3.0 ; $ws.Range("M8").Value = 4.752338666666668 ; $ws.Range("N8").Value = 14.257016 ; $ws.Range("O8").Value = 0.04151639666945049 ; $ws.Range("P8").Value = 0.04151639666945049 ; $ws.Range("Q8").Value = 7238.808433957379 ; $ws.Range("R8").Value = 65149.27590561641 ; $ws.Range("S8").Value = 0.03547091209410772 ; $ws.Range("T8").Value = 0.03547091209410772
$ws.Range("E9").Value = 3.0 ; $ws.Range("G9").Value = 1523.209716666667 ; $ws.Range("H9").Value = 4569.62915 ; $ws.Range("I9").Value = 0.854383206146807 ; $ws.Range("J9").Value = 0.854383206146807 ; $ws.Range("K9").Value = 3.0 ; $ws.Range("M9").Value = 7.095953333333333 ; $ws.Range("N9").Value = 21.28786 ; $ws.Range("O9").Value = 0.06199019766855336 ; $ws.Range("P9").Value = 0.06199019766855336 ; $ws.Range("Q9").Value = 10808.62506634655 ; $ws.Range("R9").Value = 97277.62559711898 ; $ws.Range("S9").Value = 0.05296338383373294 ; $ws.Range("T9").Value = 0.05296338383373294
$ws.Range("E10").Value = 3.0 ; $ws.Range("G10").Value = 1523.209716666667 ; $ws.Range("H10").Value = 4569.62915 ; $ws.Range("I10").Value = 0.854383206146807 ; $ws.Range("J10").Value = 0.854383206146807 ; $ws.Range("K10").Value = 3.0 ; $ws.Range("M10").Value = 27.34521433333333 ; $ws.Range("N10").Value = 82.03564300000001 ; $ws.Range("O10").Value = 0.238887597223811 ; $ws.Range("P10").Value = 0.238887597223811 ; $ws.Range("Q10").Value = 41652.49617686593 ; $ws.Range("R10").Value = 374872.4655917935 ; $ws.Range("S10").Value = 0.2041015512247867 ; $ws.Range("T10").Value = 0.2041015512247867
$ws.Range("E11").Value = 3.0 ; $ws.Range("G11").Value = 1523.209716666667 ; $ws.Range("H11").Value = 4569.62915 ; $ws.Range("I11").Value = 0.854383206146807 ; $ws.Range("J11").Value = 0.854383206146807 ; $ws.Range("K11").Value = 3.0 ; $ws.Range("M11").Value = 19.73820233333333 ; $ws.Range("N11").Value = 59.214607 ; $ws.Range("O11").Value = 0.1724327946912327 ; $ws.Range("P11").Value = 0.1724327946912327 ; $ws.Range("Q11").Value = 30065.421583666 ; $ws.Range("R11").Value = 270588.7942529941 ; $ws.Range("S11").Value = 0.1473236839731495 ; $ws.Range("T11").Value = 0.1473236839731495
$ws.Range("E12").Value = 3.0 ; $ws.Range("G12").Value = 1523.209716666667 ; $ws.Range("H12").Value = 4569.62915 ; $ws.Range("I12").Value = 0.854383206146807 ; $ws.Range("J12").Value = 0.854383206146807 ; $ws.Range("K12").Value = 3.0 ; $ws.Range("M12").Value = 51.532109 ; $ws.Range("N12").Value = 154.596327 ; $ws.Range("O12").Value = 0.4501841363839443 ; $ws.Range("P12").Value = 0.4501841363839443 ; $ws.Range("Q12").Value = 78494.20914912577 ; $ws.Range("R12").Value = 706447.882342132 ; $ws.Range("S12").Value = 0.3846297658001458 ; $ws.Range("T12").Value = 0.3846297658001458
$ws.Range("E13").Value = 3.0 ; $ws.Range("G13").Value = 1523.209716666667 ; $ws.Range("H13").Value = 4569.62915 ; $ws.Range("I13").Value = 0.854383206146807 ; $ws.Range("J13").Value = 0.854383206146807 ; $ws.Range("K13").Value = 3.0 ; $ws.Range("M13").Value = 4.005140333333333 ; $ws.Range("N13").Value = 12.015421 ; $ws.Range("O13").Value = 0.03498887736300817 ; $ws.Range("P13").Value = 0.03498887736300817 ; $ws.Range("Q13").Value = 6100.668672346905 ; $ws.Range("R13").Value = 54906.01805112215 ; $ws.Range("S13").Value = 0.02989390922088436 ; $ws.Range("T13").Value = 0.02989390922088436
$ws.Range("E14").Value = 3.0 ; $ws.Range("G14").Value = 1.839165666666666 ; $ws.Range("H14").Value = 5.517497 ; $ws.Range("I14").Value = 0.001031605984211079 ; $ws.Range("J14").Value = 0.00103160598421108 ; $ws.Range("K14").Value = 3.0 ; $ws.Range("M14").Value = 4.752338666666668 ; $ws.Range("N14").Value = 14.257016 ; $ws.Range("O14").Value = 0.04151639666945049 ; $ws.Range("P14").Value = 0.04151639666945049 ; $ws.Range("Q14").Value = 8.740338112105778 ; $ws.Range("R14").Value = 78.663043008952 ; $ws.Range("S14").Value = 0.00004282856324708605 ; $ws.Range("T14").Value = 0.00004282856324708606
$ws.Range("E15").Value = 3.0 ; $ws.Range("G15").Value = 1.839165666666666 ; $ws.Range("H15").Value = 5.517497 ; $ws.Range("I15").Value = 0.001031605984211079 ; $ws.Range("J15").Value = 0.00103160598421108 ; $ws.Range("K15").Value = 3.0 ; $ws.Range("M15").Value = 7.095953333333333 ; $ws.Range("N15").Value = 21.28786 ; $ws.Range("O15").Value = 0.06199019766855336 ; $ws.Range("P15").Value = 0.06199019766855336 ; $ws.Range("Q15").Value = 13.05063374293555 ; $ws.Range("R15").Value = 117.45570368642 ; $ws.Range("S15").Value = 0.00006394945887730736 ; $ws.Range("T15").Value = 0.00006394945887730738
$ws.Range("E16").Value = 3.0 ; $ws.Range("G16").Value = 1.839165666666666 ; $ws.Range("H16").Value = 5.517497 ; $ws.Range("I16").Value = 0.001031605984211079 ; $ws.Range("J16").Value = 0.00103160598421108 ; $ws.Range("K16").Value = 3.0 ; $ws.Range("M16").Value = 27.34521433333333 ; $ws.Range("N16").Value = 82.03564300000001 ; $ws.Range("O16").Value = 0.238887597223811 ; $ws.Range("P16").Value = 0.238887597223811 ; $ws.Range("Q16").Value = 50.29237934950788 ; $ws.Range("R16").Value = 452.631414145571 ; $ws.Range("S16").Value = 0.0002464378748498895 ; $ws.Range("T16").Value = 0.0002464378748498895
$ws.Range("E17").Value = 3.0 ; $ws.Range("G17").Value = 1.839165666666666 ; $ws.Range("H17").Value = 5.517497 ; $ws.Range("I17").Value = 0.001031605984211079 ; $ws.Range("J17").Value = 0.00103160598421108 ; $ws.Range("K17").Value = 3.0 ; $ws.Range("M17").Value = 19.73820233333333 ; $ws.Range("N17").Value = 59.214607 ; $ws.Range("O17").Value = 0.1724327946912327 ; $ws.Range("P17").Value = 0.1724327946912327 ; $ws.Range("Q17").Value = 36.30182405318655 ; $ws.Range("R17").Value = 326.716416478679 ; $ws.Range("S17").Value = 0.0001778827028777161 ; $ws.Range("T17").Value = 0.0001778827028777161
$ws.Range("E18").Value = 3.0 ; $ws.Range("G18").Value = 1.839165666666666 ; $ws.Range("H18").Value = 5.517497 ; $ws.Range("I18").Value = 0.001031605984211079 ; $ws.Range("J18").Value = 0.00103160598421108 ; $ws.Range("K18").Value = 3.0 ; $ws.Range("M18").Value = 51.532109 ; $ws.Range("N18").Value = 154.596327 ; $ws.Range("O18").Value = 0.4501841363839443 ; $ws.Range("P18").Value = 0.4501841363839443 ; $ws.Range("Q18").Value = 94.77608560372433 ; $ws.Range("R18").Value = 852.984770433519 ; $ws.Range("S18").Value = 0.0004644126490905738 ; $ws.Range("T18").Value = 0.0004644126490905738
$ws.Range("E19").Value = 3.0 ; $ws.Range("G19").Value = 1.839165666666666 ; $ws.Range("H19").Value = 5.517497 ; $ws.Range("I19").Value = 0.001031605984211079 ; $ws.Range("J19").Value = 0.00103160598421108 ; $ws.Range("K19").Value = 3.0 ; $ws.Range("M19").Value = 4.005140333333333 ; $ws.Range("N19").Value = 12.015421 ; $ws.Range("O19").Value = 0.03498887736300817 ; $ws.Range("P19").Value = 0.03498887736300817 ; $ws.Range("Q19").Value = 7.366116591248555 ; $ws.Range("R19").Value = 66.295049321237 ; $ws.Range("S19").Value = 0.00003609473526850681 ; $ws.Range("T19").Value = 0.00003609473526850682
$ws.Range("E20").Value = 3.0 ; $ws.Range("G20").Value = 8.497059666666667 ; $ws.Range("H20").Value = 25.491179 ; $ws.Range("I20").Value = 0.004766083751562676 ; $ws.Range("J20").Value = 0.004766083751562675 ; $ws.Range("K20").Value = 3.0 ; $ws.Range("M20").Value = 4.752338666666668 ; $ws.Range("N20").Value = 14.257016 ; $ws.Range("O20").Value = 0.04151639666945049 ; $ws.Range("P20").Value = 0.04151639666945049 ; $ws.Range("Q20").Value = 40.38090520687378 ; $ws.Range("R20").Value = 363.428146861864 ; $ws.Range("S20").Value = 0.0001978706235896987 ; $ws.Range("T20").Value = 0.0001978706235896987
$ws.Range("E21").Value = 3.0 ; $ws.Range("G21").Value = 8.497059666666667 ; $ws.Range("H21").Value = 25.491179 ; $ws.Range("I21").Value = 0.004766083751562676 ; $ws.Range("J21").Value = 0.004766083751562675 ; $ws.Range("K21").Value = 3.0 ; $ws.Range("M21").Value = 7.095953333333333 ; $ws.Range("N21").Value = 21.28786 ; $ws.Range("O21").Value = 0.06199019766855336 ; $ws.Range("P21").Value = 0.06199019766855336 ; $ws.Range("Q21").Value = 60.29473886521556 ; $ws.Range("R21").Value = 542.65264978694 ; $ws.Range("S21").Value = 0.0002954504738642506 ; $ws.Range("T21").Value = 0.0002954504738642506
$ws.Range("E22").Value = 3.0 ; $ws.Range("G22").Value = 8.497059666666667 ; $ws.Range("H22").Value = 25.491179 ; $ws.Range("I22").Value = 0.004766083751562676 ; $ws.Range("J22").Value = 0.004766083751562675 ; $ws.Range("K22").Value = 3.0 ; $ws.Range("M22").Value = 27.34521433333333 ; $ws.Range("N22").Value = 82.03564300000001 ; $ws.Range("O22").Value = 0.238887597223811 ; $ws.Range("P22").Value = 0.238887597223811 ; $ws.Range("Q22").Value = 232.3539177881219 ; $ws.Range("R22").Value = 2091.185260093097 ; $ws.Range("S22").Value = 0.001138558295578254 ; $ws.Range("T22").Value = 0.001138558295578254
$ws.Range("E23").Value = 3.0 ; $ws.Range("G23").Value = 8.497059666666667 ; $ws.Range("H23").Value = 25.491179 ; $ws.Range("I23").Value = 0.004766083751562676 ; $ws.Range("J23").Value = 0.004766083751562675 ; $ws.Range("K23").Value = 3.0 ; $ws.Range("M23").Value = 19.73820233333333 ; $ws.Range("N23").Value = 59.214607 ; $ws.Range("O23").Value = 0.1724327946912327 ; $ws.Range("P23").Value = 0.1724327946912327 ; $ws.Range("Q23").Value = 167.7166829390726 ; $ws.Range("R23").Value = 1509.450146451653 ; $ws.Range("S23").Value = 0.0008218291410144267 ; $ws.Range("T23").Value = 0.0008218291410144266
$ws.Range("E24").Value = 3.0 ; $ws.Range("G24").Value = 8.497059666666667 ; $ws.Range("H24").Value = 25.491179 ; $ws.Range("I24").Value = 0.004766083751562676 ; $ws.Range("J24").Value = 0.004766083751562675 ; $ws.Range("K24").Value = 3.0 ; $ws.Range("M24").Value = 51.532109 ; $ws.Range("N24").Value = 154.596327 ; $ws.Range("O24").Value = 0.4501841363839443 ; $ws.Range("P24").Value = 0.4501841363839443 ; $ws.Range("Q24").Value = 437.8714049221703 ; $ws.Range("R24").Value = 3940.842644299533 ; $ws.Range("S24").Value = 0.002145615297630793 ; $ws.Range("T24").Value = 0.002145615297630792
$ws.Range("E25").Value = 3.0 ; $ws.Range("G25").Value = 8.497059666666667 ; $ws.Range("H25").Value = 25.491179 ; $ws.Range("I25").Value = 0.004766083751562676 ; $ws.Range("J25").Value = 0.004766083751562675 ; $ws.Range("K25").Value = 3.0 ; $ws.Range("M25").Value = 4.005140333333333 ; $ws.Range("N25").Value = 12.015421 ; $ws.Range("O25").Value = 0.03498887736300817 ; $ws.Range("P25").Value = 0.03498887736300817 ; $ws.Range("Q25").Value = 34.03191638570656 ; $ws.Range("R25").Value = 306.287247471359 ; $ws.Range("S25").Value = 0.0001667599198852524 ; $ws.Range("T25").Value = 0.0001667599198852524
$ws.Range("E26").Value = 3.0 ; $ws.Range("G26").Value = 1.619322333333334 ; $ws.Range("H26").Value = 4.857967 ; $ws.Range("I26").Value = 0.0009082937114963446 ; $ws.Range("J26").Value = 0.0009082937114963446 ; $ws.Range("K26").Value = 3.0 ; $ws.Range("M26").Value = 4.752338666666668 ; $ws.Range("N26").Value = 14.257016 ; $ws.Range("O26").Value = 0.04151639666945049 ; $ws.Range("P26").Value = 0.04151639666945049 ; $ws.Range("Q26").Value = 7.695568138496891 ; $ws.Range("R26").Value = 69.26011324647202 ; $ws.Range("S26").Value = 0.00003770908201884966 ; $ws.Range("T26").Value = 0.00003770908201884966
$ws.Range("E27").Value = 3.0 ; $ws.Range("G27").Value = 1.619322333333334 ; $ws.Range("H27").Value = 4.857967 ; $ws.Range("I27").Value = 0.0009082937114963446 ; $ws.Range("J27").Value = 0.0009082937114963446 ; $ws.Range("K27").Value = 3.0 ; $ws.Range("M27").Value = 7.095953333333333 ; $ws.Range("N27").Value = 21.28786 ; $ws.Range("O27").Value = 0.06199019766855336 ; $ws.Range("P27").Value = 0.06199019766855336 ; $ws.Range("Q27").Value = 11.49063570895778 ; $ws.Range("R27").Value = 103.41572138062 ; $ws.Range("S27").Value = 0.00005630530671676238 ; $ws.Range("T27").Value = 0.00005630530671676238
$ws.Range("E28").Value = 3.0 ; $ws.Range("G28").Value = 1.619322333333334 ; $ws.Range("H28").Value = 4.857967 ; $ws.Range("I28").Value = 0.0009082937114963446 ; $ws.Range("J28").Value = 0.0009082937114963446 ; $ws.Range("K28").Value = 3.0 ; $ws.Range("M28").Value = 27.34521433333333 ; $ws.Range("N28").Value = 82.03564300000001 ; $ws.Range("O28").Value = 0.238887597223811 ; $ws.Range("P28").Value = 0.238887597223811 ; $ws.Range("Q28").Value = 44.28071627975345 ; $ws.Range("R28").Value = 398.526446517781 ; $ws.Range("S28").Value = 0.0002169801023128591 ; $ws.Range("T28").Value = 0.0002169801023128591
$ws.Range("E29").Value = 3.0 ; $ws.Range("G29").Value = 1.619322333333334 ; $ws.Range("H29").Value = 4.857967 ; $ws.Range("I29").Value = 0.0009082937114963446 ; $ws.Range("J29").Value = 0.0009082937114963446 ; $ws.Range("K29").Value = 3.0 ; $ws.Range("M29").Value = 19.73820233333333 ; $ws.Range("N29").Value = 59.214607 ; $ws.Range("O29").Value = 0.1724327946912327 ; $ws.Range("P29").Value = 0.1724327946912327 ; $ws.Range("Q29").Value = 31.96251185821878 ; $ws.Range("R29").Value = 287.662606723969 ; $ws.Range("S29").Value = 0.0001566196230737869 ; $ws.Range("T29").Value = 0.0001566196230737869
$ws.Range("E30").Value = 3.0 ; $ws.Range("G30").Value = 1.619322333333334 ; $ws.Range("H30").Value = 4.857967 ; $ws.Range("I30").Value = 0.0009082937114963446 ; $ws.Range("J30").Value = 0.0009082937114963446 ; $ws.Range("K30").Value = 3.0 ; $ws.Range("M30").Value = 51.532109 ; $ws.Range("N30").Value = 154.596327 ; $ws.Range("O30").Value = 0.4501841363839443 ; $ws.Range("P30").Value = 0.4501841363839443 ; $ws.Range("Q30").Value = 83.44709498746768 ; $ws.Range("R30").Value = 751.023854887209 ; $ws.Range("S30").Value = 0.0004088994200929494 ; $ws.Range("T30").Value = 0.0004088994200929494
$ws.Range("E31").Value = 3.0 ; $ws.Range("G31").Value = 1.619322333333334 ; $ws.Range("H31").Value = 4.857967 ; $ws.Range("I31").Value = 0.0009082937114963446 ; $ws.Range("J31").Value = 0.0009082937114963446 ; $ws.Range("K31").Value = 3.0 ; $ws.Range("M31").Value = 4.005140333333333 ; $ws.Range("N31").Value = 12.015421 ; $ws.Range("O31").Value = 0.03498887736300817 ; $ws.Range("P31").Value = 0.03498887736300817 ; $ws.Range("Q31").Value = 6.485613189900779 ; $ws.Range("R31").Value = 58.370518709107 ; $ws.Range("S31").Value = 0.00003178017728113713 ; $ws.Range("T31").Value = 0.00003178017728113713
$ws.Range("E32").Value = 3.0 ; $ws.Range("G32").Value = 196.5440926666667 ; $ws.Range("H32").Value = 589.632278 ; $ws.Range("I32").Value = 0.1102435010782627 ; $ws.Range("J32").Value = 0.1102435010782627 ; $ws.Range("K32").Value = 3.0 ; $ws.Range("M32").Value = 4.752338666666668 ; $ws.Range("N32").Value = 14.257016 ; $ws.Range("O32").Value = 0.04151639666945049 ; $ws.Range("P32").Value = 0.04151639666945049 ; $ws.Range("Q32").Value = 934.0440912847166 ; $ws.Range("R32").Value = 8406.39682156245 ; $ws.Range("S32").Value = 0.004576912920994144 ; $ws.Range("T32").Value = 0.004576912920994144
$ws.Range("E33").Value = 3.0 ; $ws.Range("G33").Value = 196.5440926666667 ; $ws.Range("H33").Value = 589.632278 ; $ws.Range("I33").Value = 0.1102435010782627 ; $ws.Range("J33").Value = 0.1102435010782627 ; $ws.Range("K33").Value = 3.0 ; $ws.Range("M33").Value = 7.095953333333333 ; $ws.Range("N33").Value = 21.28786 ; $ws.Range("O33").Value = 0.06199019766855336 ; $ws.Range("P33").Value = 0.06199019766855336 ; $ws.Range("Q33").Value = 1394.667709505009 ; $ws.Range("R33").Value = 12552.00938554508 ; $ws.Range("S33").Value = 0.006834016423514878 ; $ws.Range("T33").Value = 0.006834016423514878
$ws.Range("E34").Value = 3.0 ; $ws.Range("G34").Value = 196.5440926666667 ; $ws.Range("H34").Value = 589.632278 ; $ws.Range("I34").Value = 0.1102435010782627 ; $ws.Range("J34").Value = 0.1102435010782627 ; $ws.Range("K34").Value = 3.0 ; $ws.Range("M34").Value = 27.34521433333333 ; $ws.Range("N34").Value = 82.03564300000001 ; $ws.Range("O34").Value = 0.238887597223811 ; $ws.Range("P34").Value = 0.238887597223811 ; $ws.Range("Q34").Value = 5374.540339920529 ; $ws.Range("R34").Value = 48370.86305928476 ; $ws.Range("S34").Value = 0.02633580508212678 ; $ws.Range("T34").Value = 0.02633580508212678
$ws.Range("E35").Value = 3.0 ; $ws.Range("G35").Value = 196.5440926666667 ; $ws.Range("H35").Value = 589.632278 ; $ws.Range("I35").Value = 0.1102435010782627 ; $ws.Range("J35").Value = 0.1102435010782627 ; $ws.Range("K35").Value = 3.0 ; $ws.Range("M35").Value = 19.73820233333333 ; $ws.Range("N35").Value = 59.214607 ; $ws.Range("O35").Value = 0.1724327946912327 ; $ws.Range("P35").Value = 0.1724327946912327 ; $ws.Range("Q35").Value = 3879.427068476083 ; $ws.Range("R35").Value = 34914.84361628475 ; $ws.Range("S35").Value = 0.01900959498747075 ; $ws.Range("T35").Value = 0.01900959498747075
$ws.Range("E36").Value = 3.0 ; $ws.Range("G36").Value = 196.5440926666667 ; $ws.Range("H36").Value = 589.632278 ; $ws.Range("I36").Value = 0.1102435010782627 ; $ws.Range("J36").Value = 0.1102435010782627 ; $ws.Range("K36").Value = 3.0 ; $ws.Range("M36").Value = 51.532109 ; $ws.Range("N36").Value = 154.596327 ; $ws.Range("O36").Value = 0.4501841363839443 ; $ws.Range("P36").Value = 0.4501841363839443 ; $ws.Range("Q36").Value = 10128.33160660477 ; $ws.Range("R36").Value = 91154.98445944292 ; $ws.Range("S36").Value = 0.04962987532486011 ; $ws.Range("T36").Value = 0.04962987532486011
$ws.Range("E37").Value = 3.0 ; $ws.Range("G37").Value = 196.5440926666667 ; $ws.Range("H37").Value = 589.632278 ; $ws.Range("I37").Value = 0.1102435010782627 ; $ws.Range("J37").Value = 0.1102435010782627 ; $ws.Range("K37").Value = 3.0 ; $ws.Range("M37").Value = 4.005140333333333 ; $ws.Range("N37").Value = 12.015421 ; $ws.Range("O37").Value = 0.03498887736300817 ; $ws.Range("P37").Value = 0.03498887736300817 ; $ws.Range("Q37").Value = 787.1866728176709 ; $ws.Range("R37").Value = 7084.680055359038 ; $ws.Range("S37").Value = 0.003857296339295992 ; $ws.Range("T37").Value = 0.003857296339295992
